$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (12-21) are the existing match-log rows (2-11) re-ordered and
# appended beneath the current table. Force text formatting first so the
# numeric-looking values (runs/balls/4s/6s/strike-rate) are written as text,
# matching the rest of the sheet instead of being coerced to numbers.
$ws.Range("A12:K21").NumberFormat = "@"

$data = @(
  @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Sunrisers Hyderabad", "Chennai Super Kings", "Jonny Bairstow †", "23", "24", "2", "0", "95.83"),
  @(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Sunrisers Hyderabad", "Rajasthan Royals", "Jonny Bairstow †", "10", "7", "1", "1", "142.85"),
  @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Sunrisers Hyderabad", "Kolkata Knight Riders", "Jonny Bairstow †", "36", "28", "7", "0", "128.57"),
  @(" Dubai (DSC)", " October 02 2020", "Sunrisers won by 7 runs", "Sunrisers Hyderabad", "Chennai Super Kings", "Jonny Bairstow †", "0", "3", "0", "0", "0.00"),
  @(" Abu Dhabi", " September 29 2020", "Sunrisers won by 15 runs", "Sunrisers Hyderabad", "Delhi Capitals", "Jonny Bairstow †", "53", "48", "2", "1", "110.41"),
  @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Sunrisers Hyderabad", "Mumbai Indians", "Jonny Bairstow †", "25", "15", "2", "2", "166.66"),
  @(" Dubai (DSC)", " September 21 2020", "RCB won by 10 runs", "Sunrisers Hyderabad", "Royal Challengers Bangalore", "Jonny Bairstow †", "61", "43", "6", "2", "141.86"),
  @(" Dubai (DSC)", " October 11 2020", "Royals won by 5 wickets (with 1 ball remaining)", "Sunrisers Hyderabad", "Rajasthan Royals", "Jonny Bairstow †", "16", "19", "0", "1", "84.21"),
  @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Sunrisers Hyderabad", "Kings XI Punjab", "Jonny Bairstow †", "19", "20", "4", "0", "95.00"),
  @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Sunrisers Hyderabad", "Kings XI Punjab", "Jonny Bairstow †", "97", "55", "7", "6", "176.36")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = 12 + $i
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, 1 + $j).Value = $row[$j]
  }
}

# Extend the "number stored as text" ignored-error suppression (the green
# triangle indicator) over the freshly written numeric-looking text cells,
# same as the rest of the table (A1:K11 -> A1:K21).
try {
  $ws.Range("A1:K21").Errors.Item(9).Ignore = $true
} catch {
}
